$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsProfit = $wb.Worksheets.Item("Profit")
$wsConstraints = $wb.Worksheets.Item("Constraints")

# Remove the now-redundant "Constraints" sheet; its data is folded into "Profit" (renamed "Data").
[void]$wsConstraints.Delete()

# Rename the remaining sheet to "Data".
$wsProfit.Name = "Data"
$ws = $wsProfit

# Header row gains the "Available" column.
$ws.Range("D1").Value = "Available"

# Shift the Profit data row down to make room, and add the two constraint rows below it.
$ws.Range("A3").Value = "F1"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 4000

$ws.Range("A4").Value = "F2"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 5000

# Apply the numeric format used throughout the data block.
$ws.Range("B2:D4").NumberFormat = "0.00"

[void]$ws.Range("F10").Select()
